# Add final summary slide with final stats, and annotate the two
# "Some other interesting stats. . . " slides with their distances.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Update the two "Some other interesting stats. . . " title slides
#    with the distance annotations.
# ---------------------------------------------------------------------
$s5k = $p.Slides.Item(24)
$s5k.Shapes.Item(1).TextFrame.TextRange.Text = "Some other interesting stats. . . (5km) "

$s10k = $p.Slides.Item(36)
$s10k.Shapes.Item(1).TextFrame.TextRange.Text = "Some other interesting stats. . . (10km) "

# ---------------------------------------------------------------------
# 2) Add the new final "Finally. . . " summary slide at the end of the
#    deck, using the same "Title and Content" layout as the preceding
#    "Conclusions" slide.
# ---------------------------------------------------------------------
$newIndex = $p.Slides.Count + 1
$slide = $p.Slides.Add($newIndex, 2)

# --- Title ---
$title = $slide.Shapes.Item(1)
$ttr = $title.TextFrame.TextRange
$ttr.Text = "Finally. . . "
$ttr.LanguageID = "en-GB"

# --- Body content placeholder ---
$body = $slide.Shapes.Item(2)
$tr = $body.TextFrame.TextRange
$dash = [char]0x2013

# First paragraph, built up run by run so the segment boundaries match
# the authored runs. Setting LanguageID while the range is still a
# single run lets every subsequently-inserted run inherit "en-GB".
$tr.Text = "You ran 1379.65 kms for a total time of 171:35:33 ("
$tr.LanguageID = "en-GB"

[void]$tr.InsertAfter("hh:mm:ss")
[void]$tr.InsertAfter(") (over 7 days of running!!), over 603 days, burning ")
[void]$tr.InsertAfter("146646 calories ")
[void]$tr.InsertAfter("between the 03/12/12 $dash 30/07/19. ")

# Blank paragraph.
[void]$tr.InsertAfter("`r")
[void]$tr.InsertAfter("X`r")

# Final paragraph: "WELL DONE! " + Wingdings smiley.
[void]$tr.InsertAfter("WELL DONE! ")
[void]$tr.InsertAfter("J")

# Remove the "X" placeholder so the middle paragraph is blank.
$xIndex = $tr.Text.IndexOf("X")
$xchar = $tr.Characters($xIndex + 1, 1)
$xchar.Text = ""

# Apply the Wingdings smiley font to the trailing "J" glyph.
$sym = $tr.Characters($tr.Length, 1)
$sym.Font.Name = "Wingdings"
